$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.663.68"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.157.64"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "616.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.153.40"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.473"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("D15").Value = "3.680.88"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").Value = "64.634.19"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "3.160.30"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.720"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -3.65%  "
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("E29").Value = "  -6.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("E31").Value = "  -8.28%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("D36").Value = "0.0₃0779"
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "460.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  -3.70%  "
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").Value = "2.845.93"
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.14%  "
